$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new work item ("Short point (up to 3 mtr.)") was inserted as the new row 9, pushing
# the existing item rows down by one. Reproduce that by inserting a row just above the
# current "Grand Total" summary block (old row 18) which mirrors the row that opens up
# once the bottom block is shifted, then fill in every cell with the final values.
$ws.Rows(18).Insert()

# --- Row 8: quantity changes, everything else the same ---
$ws.Range("C8").Value = 90

# --- Row 9: brand-new line item ---
$ws.Range("A9").Value = "P. point"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 86
$ws.Range("D9").Value = "'2"
$ws.Range("E9").Value = "Short point (up to 3 mtr.)"
$ws.Range("F9").Value = 256
$ws.Range("G9").Value = "'22016.00"
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = ""

# --- Row 10 ---
$ws.Range("C10").Value = 11
$ws.Range("D10").Value = "'3"
$ws.Range("E10").Value = "Medium point (up to 6 mtr.)"
$ws.Range("F10").Value = 472
$ws.Range("G10").Value = "'5192.00"

# --- Row 11 ---
$ws.Range("A11").Value = "P. point"
$ws.Range("C11").Value = 22
$ws.Range("D11").Value = "'4"
$ws.Range("E11").Value = "Long point  (up to 10 mtr.)"
$ws.Range("F11").Value = 662
$ws.Range("G11").Value = "'14564.00"

# --- Row 12 ---
$ws.Range("A12").Value = ""
$ws.Range("C12").Value = 93
$ws.Range("D12").Value = "'2.0"
$ws.Range("E12").Value = "Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it's  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet's & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = "'0.00"

# --- Row 13 ---
$ws.Range("A13").Value = "P. point"
$ws.Range("C13").Value = 56
$ws.Range("D13").Value = "'6"
$ws.Range("E13").Value = "On board"
$ws.Range("F13").Value = 136
$ws.Range("G13").Value = "'7616.00"

# --- Row 14 ---
$ws.Range("A14").Value = "Each"
$ws.Range("C14").Value = 69
$ws.Range("D14").Value = "'3.0"
$ws.Range("E14").Value = 'P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F14").Value = 23
$ws.Range("G14").Value = "'1587.00"

# --- Row 15 ---
$ws.Range("A15").Value = ""
$ws.Range("C15").Value = 55
$ws.Range("D15").Value = "'8"
$ws.Range("E15").Value = "Total"

# --- Row 16 ---
$ws.Range("A16").Value = "%"
$ws.Range("C16").Value = 64
$ws.Range("D16").Value = "'9"
$ws.Range("E16").Value = "Add Tender Premium "

# --- Row 17: the old "Grand Total" S.No. line, now shifted here ---
$ws.Range("A17").Value = ""
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 89
$ws.Range("D17").Value = "'10"
$ws.Range("E17").Value = "Grand Total"
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = "'0.00"
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = ""

# --- Totals block (rows 19 & 21) recomputed for the new set of quantities ---
$ws.Range("G19").Value = "'50975.00"
$ws.Range("H19").Value = "'50975.00"
$ws.Range("G21").Value = "'50975.00"
$ws.Range("H21").Value = "'50975.00"
